$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35 (shifts existing rows 35-38 down to 36-39)
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).RowHeight = 13.05

# Populate the newly inserted row with the new leaderboard entry
$ws.Range("A35").Value = "Las Americas "
$ws.Range("C35").Value = "013"
$ws.Range("E35").Value = "0008279"

# Salesperson (B) and Last Invoice Date (D) are blank for this entry;
# copy the existing blank-cell style (from D34) so the number format /
# alignment matches the sheet's established blank-cell style.
$ws.Range("D34").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("D35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F35 is an empty placeholder cell like the rest of column F; touch it so
# it is materialized in the sheet even though it carries no value/style.
$ws.Range("F35").Font.Name = "Arial"
